$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.910.93'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.876.12'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7412'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -4.03%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3157'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +1.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07207'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.71'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -3.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08383'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7516'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -1.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.436'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.862.94'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.68'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.922.67'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.095'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '246.81'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.60'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007853'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9987'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.126.39'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -5.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.021'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9987'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('E25').Value = '  -4.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.271'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.19'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.66'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.040'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('E30').Value = '  +4.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.609'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +2.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.533'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.281'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +4.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05317'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.239'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7552'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +1.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.0000'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.690'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01964'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.757'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4507'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +1.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.109.59'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.067'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.62'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8560'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.16'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.632'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.856'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.480'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -3.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.928'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -2.12%  '
